$d = $word.ActiveDocument

$replacements = @(
    @("682×7=", "327×9="),
    @("436×2=", "356×7="),
    @("155×4=", "305×8="),
    @("896×6=", "367×4="),
    @("566×8=", "568×3="),
    @("715×2=", "810×3="),
    @("721×9=", "469×6="),
    @("280×7=", "856×3="),
    @("223×4=", "956×6="),
    @("532×7=", "881×3="),
    @("811×8=", "549×5="),
    @("474×5=", "714×8="),
    @("407×8=", "256×4="),
    @("287×7=", "579×5="),
    @("572×6=", "986×5="),
    @("275×3=", "370×8="),
    @("847×5=", "932×6="),
    @("304×7=", "824×7="),
    @("189×4=", "978×6="),
    @("539×5=", "933×5="),
    @("731×3=", "967×7="),
    @("455×3=", "110×7="),
    @("670×7=", "754×3="),
    @("771×4=", "569×2="),
    @("816×7=", "988×5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
